# Updated cryptos list (GitHub Actions data refresh).
# Column D = "Price", column E = "Volume(1h)" - both stored as literal
# text in the sheet (not numbers), matching the original inline-string cells.
#
# Excel auto-converts plain decimal-looking strings ("1.001", "39.26", ...)
# into real numbers when assigned through .Value (e.g. "4.800" would lose its
# trailing zero). Those values are written with a leading apostrophe, Excel's
# standard "force text" marker, which is stripped from the stored value but
# keeps the cell as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.956.92"
$ws.Range("E2").Value = "  -7.94%  "
$ws.Range("D3").Value = "1.410.78"
$ws.Range("E3").Value = "  -8.18%  "
$ws.Range("D4").Value = "'" + "1.001"   # force text, avoid numeric auto-conversion
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'" + "1.001"   # force text, avoid numeric auto-conversion
$ws.Range("D6").Value = "'" + "272.97"   # force text, avoid numeric auto-conversion
$ws.Range("E6").Value = "  -5.58%  "
$ws.Range("D7").Value = "'" + "0.3686"   # force text, avoid numeric auto-conversion
$ws.Range("E7").Value = "  -6.11%  "
$ws.Range("D8").Value = "'" + "0.3074"   # force text, avoid numeric auto-conversion
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").Value = "'" + "39.26"   # force text, avoid numeric auto-conversion
$ws.Range("E9").Value = "  -7.54%  "
$ws.Range("D10").Value = "'" + "0.9968"   # force text, avoid numeric auto-conversion
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").Value = "'" + "0.06553"   # force text, avoid numeric auto-conversion
$ws.Range("E11").Value = "  -8.61%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'" + "5.331"   # force text, avoid numeric auto-conversion
$ws.Range("E13").Value = "  -4.97%  "
$ws.Range("D14").Value = "'" + "16.94"   # force text, avoid numeric auto-conversion
$ws.Range("E14").Value = "  -8.34%  "
$ws.Range("D15").Value = "'" + "6.147"   # force text, avoid numeric auto-conversion
$ws.Range("E15").Value = "  -6.85%  "
$ws.Range("D16").Value = "1.411.14"
$ws.Range("E16").Value = "  -8.32%  "
$ws.Range("D17").Value = "'" + "0.00001008"   # force text, avoid numeric auto-conversion
$ws.Range("E17").Value = "  -7.94%  "
$ws.Range("D18").Value = "'" + "0.05742"   # force text, avoid numeric auto-conversion
$ws.Range("E18").Value = "  -12.77%  "
$ws.Range("D19").Value = "'" + "73.32"   # force text, avoid numeric auto-conversion
$ws.Range("E19").Value = "  -11.64%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'" + "5.596"   # force text, avoid numeric auto-conversion
$ws.Range("E21").Value = "  -8.54%  "
$ws.Range("D22").Value = "'" + "14.38"   # force text, avoid numeric auto-conversion
$ws.Range("E22").Value = "  -6.31%  "
$ws.Range("D23").Value = "'" + "10.83"   # force text, avoid numeric auto-conversion
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'" + "2.284"   # force text, avoid numeric auto-conversion
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("D25").Value = "19.958.86"
$ws.Range("E25").Value = "  -7.97%  "
$ws.Range("D26").Value = "'" + "2.245"   # force text, avoid numeric auto-conversion
$ws.Range("E26").Value = "  -4.09%  "
$ws.Range("D27").Value = "'" + "138.83"   # force text, avoid numeric auto-conversion
$ws.Range("E27").Value = "  -5.44%  "
$ws.Range("D28").Value = "'" + "16.78"   # force text, avoid numeric auto-conversion
$ws.Range("E28").Value = "  -8.41%  "
$ws.Range("D29").Value = "1.569.09"
$ws.Range("E29").Value = "  -8.36%  "
$ws.Range("D30").Value = "'" + "108.84"   # force text, avoid numeric auto-conversion
$ws.Range("E30").Value = "  -6.98%  "
$ws.Range("D31").Value = "'" + "3.797"   # force text, avoid numeric auto-conversion
$ws.Range("E31").Value = "  -21.50%  "
$ws.Range("D32").Value = "'" + "5.322"   # force text, avoid numeric auto-conversion
$ws.Range("E32").Value = "  -9.49%  "
$ws.Range("D33").Value = "'" + "0.8265"   # force text, avoid numeric auto-conversion
$ws.Range("E33").Value = "  -13.83%  "
$ws.Range("D34").Value = "'" + "0.07680"   # force text, avoid numeric auto-conversion
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("D35").Value = "'" + "8.436"   # force text, avoid numeric auto-conversion
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").Value = "'" + "0.05754"   # force text, avoid numeric auto-conversion
$ws.Range("E36").Value = "  -5.04%  "
$ws.Range("D37").Value = "'" + "4.800"   # force text, avoid numeric auto-conversion
$ws.Range("E37").Value = "  -5.81%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'" + "0.1939"   # force text, avoid numeric auto-conversion
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("D40").Value = "'" + "0.02037"   # force text, avoid numeric auto-conversion
$ws.Range("E40").Value = "  -7.00%  "
$ws.Range("D41").Value = "'" + "10.25"   # force text, avoid numeric auto-conversion
$ws.Range("E41").Value = "  -4.02%  "
$ws.Range("D42").Value = "'" + "1.059"   # force text, avoid numeric auto-conversion
$ws.Range("E42").Value = "  -10.28%  "
$ws.Range("D43").Value = "'" + "1.276"   # force text, avoid numeric auto-conversion
$ws.Range("E43").Value = "  -11.00%  "
$ws.Range("D44").Value = "'" + "0.5281"   # force text, avoid numeric auto-conversion
$ws.Range("E44").Value = "  -7.59%  "
$ws.Range("D45").Value = "'" + "3.522"   # force text, avoid numeric auto-conversion
$ws.Range("E45").Value = "  -5.60%  "
$ws.Range("D46").Value = "'" + "12.11"   # force text, avoid numeric auto-conversion
$ws.Range("E46").Value = "  -6.68%  "
$ws.Range("D47").Value = "'" + "0.5098"   # force text, avoid numeric auto-conversion
$ws.Range("E47").Value = "  -6.71%  "
$ws.Range("D48").Value = "'" + "1.799"   # force text, avoid numeric auto-conversion
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "'" + "111.33"   # force text, avoid numeric auto-conversion
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("D50").Value = "'" + "1.038"   # force text, avoid numeric auto-conversion
$ws.Range("E50").Value = "  -10.71%  "
$ws.Range("E51").Value = "  -0.05%  "
